# Insert a new data row at row 847 ("2026/02/20", "金", 16, 201).
# All rows from the old 847 onward shift down by one (old 888 -> new 889).
#
# We copy row 846 (same date "2026/02/20" / weekday "金") and insert it at
# row 847 so the new row inherits the existing plain-text cell formatting
# (no style) instead of Excel's COM layer auto-detecting the "yyyy/mm/dd"
# string as a date literal and stamping a date number format on it.
# Only the "time" value (column C) then needs to change from 13 to 16;
# column D already matches the target (201).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(846).Copy()
$ws.Rows.Item(847).Insert()

$ws.Range("C847").Value = 16
